$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as plain text, avoiding Excel auto-converting
# numeric-looking strings into numbers (which would drop formatting like
# trailing zeros, e.g. "146.00" -> 146) and avoiding leftover custom
# number-format styles by resetting the cell style afterwards.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.723.60"
$ws.Range("E2").Value = "  +1.90%  "

$ws.Range("D3").Value = "1.623.66"
$ws.Range("E3").Value = "  +2.23%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.15%  "

Set-TextValue $ws.Range("D5") "214.57"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("E9").Value = "  +0.70%  "

Set-TextValue $ws.Range("D10") "19.35"
$ws.Range("E10").Value = "  +0.33%  "

Set-TextValue $ws.Range("D11") "0.0859"
$ws.Range("E11").Value = "  +1.11%  "

$ws.Range("D12").Value = "1.851.36"
$ws.Range("E12").Value = "  +2.17%  "

$ws.Range("D13").Value = "1.658.66"
$ws.Range("E13").Value = "  +3.74%  "

$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("E16").Value = "  -1.41%  "

$ws.Range("D17").Value = "26.708.37"
$ws.Range("E17").Value = "  +1.78%  "

Set-TextValue $ws.Range("D18") "230.76"
$ws.Range("E18").Value = "  +8.42%  "

Set-TextValue $ws.Range("D19") "7.74"
$ws.Range("E19").Value = "  +4.86%  "

$ws.Range("D20").Value = "0.0₃0728"
$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("E22").Value = "  +2.84%  "

$ws.Range("E23").Value = "  +2.70%  "

Set-TextValue $ws.Range("D24") "9.12"
$ws.Range("E24").Value = "  +1.22%  "

Set-TextValue $ws.Range("D25") "146.00"
$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +2.40%  "

Set-TextValue $ws.Range("D29") "15.60"
$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D33").Value = "1.467.27"
$ws.Range("E33").Value = "  +9.25%  "

$ws.Range("E34").Value = "  +2.51%  "

$ws.Range("E35").Value = "  -0.82%  "

$ws.Range("E36").Value = "  +1.41%  "

Set-TextValue $ws.Range("D37") "0.575"
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("E38").Value = "  +0.18%  "

Set-TextValue $ws.Range("D39") "0.841"

Set-TextValue $ws.Range("D40") "5.95"
$ws.Range("E40").Value = "  +3.72%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +2.84%  "

Set-TextValue $ws.Range("D43") "0.954"
$ws.Range("E43").Value = "  -3.53%  "

$ws.Range("D44").Value = "1.762.36"
$ws.Range("E44").Value = "  +2.21%  "

Set-TextValue $ws.Range("D45") "0.765"
$ws.Range("E45").Value = "  -0.39%  "

Set-TextValue $ws.Range("D46") "62.12"
$ws.Range("E46").Value = "  +1.35%  "

Set-TextValue $ws.Range("D47") "88.13"
$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("E48").Value = "  +2.25%  "

Set-TextValue $ws.Range("D50") "0.0964"
$ws.Range("E50").Value = "  -1.18%  "

Set-TextValue $ws.Range("D51") "7.48"
$ws.Range("E51").Value = "  +2.10%  "
